# Auto-generated edit script applying diff changes to Sheet1
# Updates odds values in rows 2, 5, 6, 8, 9, 11, 14 (columns G..BD)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 4.5
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 1.67
$ws.Range("J2").Value = 4.5
$ws.Range("K2").Value = 2.5
$ws.Range("L2").Value = 2.2
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 5.5
$ws.Range("W2").Value = 17
$ws.Range("X2").Value = 26
$ws.Range("Y2").Value = 15
$ws.Range("Z2").Value = 51
$ws.Range("AA2").Value = 34
$ws.Range("AB2").Value = 34
$ws.Range("AF2").Value = 41
$ws.Range("AI2").Value = 10
$ws.Range("AK2").Value = 15
$ws.Range("AL2").Value = 12
$ws.Range("AM2").Value = 21
$ws.Range("AN2").Value = 6.5
$ws.Range("AO2").Value = 21
$ws.Range("AP2").Value = 26
$ws.Range("AQ2").Value = 67
$ws.Range("AR2").Value = 81
$ws.Range("AS2").Value = 151
$ws.Range("AU2").Value = 7.5
$ws.Range("AX2").Value = 4
$ws.Range("AY2").Value = 8.5
$ws.Range("BA2").Value = 26
$ws.Range("BC2").Value = 101

# Row 5
$ws.Range("G5").Value = 1.65
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 5.75
$ws.Range("J5").Value = 2.38
$ws.Range("L5").Value = 7
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("W5").Value = 4.75
$ws.Range("X5").Value = 6
$ws.Range("Z5").Value = 12
$ws.Range("AA5").Value = 17
$ws.Range("AD5").Value = 7.5
$ws.Range("AE5").Value = 26
$ws.Range("AH5").Value = 10
$ws.Range("AI5").Value = 26
$ws.Range("AJ5").Value = 21
$ws.Range("AM5").Value = 67
$ws.Range("BB5").Value = 251

# Row 6
$ws.Range("G6").Value = 1.95
$ws.Range("I6").Value = 3.9
$ws.Range("J6").Value = 2.6
$ws.Range("L6").Value = 4.33
$ws.Range("W6").Value = 7.5
$ws.Range("X6").Value = 9.5
$ws.Range("Y6").Value = 8.5
$ws.Range("Z6").Value = 17
$ws.Range("AA6").Value = 15
$ws.Range("AE6").Value = 15
$ws.Range("AH6").Value = 12
$ws.Range("AI6").Value = 21
$ws.Range("AO6").Value = 10
$ws.Range("AQ6").Value = 34
$ws.Range("AW6").Value = 126
$ws.Range("AY6").Value = 21
$ws.Range("AZ6").Value = 29
$ws.Range("BD6").Value = 126

# Row 8
$ws.Range("G8").Value = 3.1
$ws.Range("I8").Value = 2.25
$ws.Range("J8").Value = 3.6
$ws.Range("U8").Value = 1.67
$ws.Range("V8").Value = 2.1
$ws.Range("AA8").Value = 23
$ws.Range("AB8").Value = 29
$ws.Range("AP8").Value = 23
$ws.Range("AR8").Value = 67

# Row 9
$ws.Range("G9").Value = 2.2
$ws.Range("H9").Value = 3.4
$ws.Range("J9").Value = 3
$ws.Range("K9").Value = 2.05
$ws.Range("S9").Value = 1.44
$ws.Range("T9").Value = 2.63
$ws.Range("U9").Value = 1.83
$ws.Range("V9").Value = 1.83
$ws.Range("AA9").Value = 19
$ws.Range("AB9").Value = 29
$ws.Range("AC9").Value = 9
$ws.Range("AG9").Value = 351
$ws.Range("AH9").Value = 8.5
$ws.Range("AP9").Value = 23
$ws.Range("AT9").Value = 2.63
$ws.Range("AV9").Value = 51
$ws.Range("BB9").Value = 81
$ws.Range("BC9").Value = 201

# Row 11
$ws.Range("G11").Value = 1.67
$ws.Range("H11").Value = 3.8
$ws.Range("I11").Value = 4.75
$ws.Range("K11").Value = 2.3
$ws.Range("M11").Value = 1.04
$ws.Range("N11").Value = 13
$ws.Range("Q11").Value = 1.8
$ws.Range("R11").Value = 2
$ws.Range("U11").Value = 1.73
$ws.Range("V11").Value = 2
$ws.Range("W11").Value = 7.5
$ws.Range("X11").Value = 8.5
$ws.Range("AB11").Value = 23
$ws.Range("AC11").Value = 12
$ws.Range("AD11").Value = 7.5
$ws.Range("AG11").Value = 201
$ws.Range("AJ11").Value = 15
$ws.Range("AN11").Value = 3.75
$ws.Range("AQ11").Value = 26
$ws.Range("AR11").Value = 41
$ws.Range("AS11").Value = 126
$ws.Range("AU11").Value = 8
$ws.Range("AZ11").Value = 29
$ws.Range("BA11").Value = 81
$ws.Range("BC11").Value = 201

# Row 14
$ws.Range("G14").Value = 2.35
$ws.Range("H14").Value = 3.2
$ws.Range("I14").Value = 2.77
$ws.Range("J14").Value = 2.9
$ws.Range("K14").Value = 2.15
$ws.Range("L14").Value = 3.35
$ws.Range("S14").Value = 1.36
$ws.Range("U14").Value = 1.55
$ws.Range("W14").Value = 10.25
$ws.Range("X14").Value = 14
$ws.Range("Y14").Value = 9
$ws.Range("Z14").Value = 27
$ws.Range("AA14").Value = 17.5
$ws.Range("AB14").Value = 22
$ws.Range("AE14").Value = 11.25
$ws.Range("AH14").Value = 10.5
$ws.Range("AI14").Value = 16
$ws.Range("AJ14").Value = 10
$ws.Range("AK14").Value = 35
$ws.Range("AL14").Value = 22
$ws.Range("AN14").Value = 4.55
$ws.Range("AO14").Value = 12.5
$ws.Range("AP14").Value = 17
$ws.Range("AQ14").Value = 50
$ws.Range("AU14").Value = 6.3
$ws.Range("AX14").Value = 5
$ws.Range("AY14").Value = 15.5
$ws.Range("BA14").Value = 70
